$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 36403
$ws.Range("J17").Value = 36403
$ws.Range("L17").Value = 109209
$ws.Range("N17").Value = -109545

$ws.Range("H31").Value = 1014279
$ws.Range("I31").Value = 1014279
$ws.Range("K31").Value = 3042837
$ws.Range("M31").Value = -3042607

$ws.Range("H43").Value = 1070.6666
$ws.Range("J43").Value = 1212
$ws.Range("L43").Value = 1212
$ws.Range("N43").Value = -1350

$ws.Range("H64").Value = 3436.7
$ws.Range("I64").Value = 3100.125
$ws.Range("J64").Value = 3559.0908
$ws.Range("K64").Value = 3100.125
$ws.Range("L64").Value = 3559.0908
$ws.Range("M64").Value = -2852.125
$ws.Range("N64").Value = -4055.0908

$ws.Range("H67").Value = 3436.7
$ws.Range("I67").Value = 3100.125
$ws.Range("J67").Value = 3559.0908
$ws.Range("K67").Value = 3100.125
$ws.Range("L67").Value = 3559.0908
$ws.Range("M67").Value = -2242.125
$ws.Range("N67").Value = -5275.0908

$ws.Range("H70").Value = 84750
$ws.Range("I70").Value = 334433.34
$ws.Range("J70").Value = 1522.2222
$ws.Range("K70").Value = 1003300.02
$ws.Range("L70").Value = 4566.6666
$ws.Range("M70").Value = -1003030.02
$ws.Range("N70").Value = -5106.6666

$ws.Range("H73").Value = 84750
$ws.Range("I73").Value = 334433.34
$ws.Range("J73").Value = 1522.2222
$ws.Range("K73").Value = 1003300.02
$ws.Range("L73").Value = 4566.6666
$ws.Range("M73").Value = -1002364.02
$ws.Range("N73").Value = -6438.6666

$ws.Range("H131").Value = 7374.56
$ws.Range("I131").Value = 1350.7142
$ws.Range("K131").Value = 4052.1426
$ws.Range("M131").Value = 987.8574000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 15254
$ws.Range("I39").Value = 13672
$ws.Range("J39").Value = 20000
$ws.Range("K39").Value = 13672
$ws.Range("L39").Value = 20000
$ws.Range("M39").Value = -13152
$ws.Range("N39").Value = -21040

$ws.Range("H45").Value = 1556.2858
$ws.Range("I45").Value = 1392.4286
$ws.Range("J45").Value = 1720.1428
$ws.Range("K45").Value = 1392.4286
$ws.Range("L45").Value = 1720.1428
$ws.Range("M45").Value = -1015.4286
$ws.Range("N45").Value = -2474.1428

$ws.Range("H74").Value = 1328.2222
$ws.Range("I74").Value = 1635.6471
$ws.Range("J74").Value = 805.6
$ws.Range("K74").Value = 1635.6471
$ws.Range("L74").Value = 805.6
$ws.Range("M74").Value = -761.6470999999999
$ws.Range("N74").Value = -2553.6

$ws.Range("H77").Value = 1328.2222
$ws.Range("I77").Value = 1635.6471
$ws.Range("J77").Value = 805.6
$ws.Range("K77").Value = 8178.2355
$ws.Range("L77").Value = 4028
$ws.Range("M77").Value = -3810.2355
$ws.Range("N77").Value = -12764

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 641.25
$ws.Range("I64").Value = 1628
$ws.Range("J64").Value = 312.33334
$ws.Range("K64").Value = 1628
$ws.Range("L64").Value = 312.33334
$ws.Range("M64").Value = -1403
$ws.Range("N64").Value = -762.33334

$ws.Range("H67").Value = 641.25
$ws.Range("I67").Value = 1628
$ws.Range("J67").Value = 312.33334
$ws.Range("K67").Value = 1628
$ws.Range("L67").Value = 312.33334
$ws.Range("M67").Value = -848
$ws.Range("N67").Value = -1872.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1495.762
$ws.Range("I16").Value = 1420.6154
$ws.Range("J16").Value = 1617.875
$ws.Range("K16").Value = 1420.6154
$ws.Range("L16").Value = 1617.875
$ws.Range("M16").Value = -1133.6154
$ws.Range("N16").Value = -2191.875

$ws.Range("H99").Value = 3616.6667
$ws.Range("J99").Value = 3850
$ws.Range("L99").Value = 3850
$ws.Range("N99").Value = -6846

$ws.Range("H113").Value = 1495.762
$ws.Range("I113").Value = 1420.6154
$ws.Range("J113").Value = 1617.875
$ws.Range("K113").Value = 1420.6154
$ws.Range("L113").Value = 1617.875
$ws.Range("M113").Value = 749.3846000000001
$ws.Range("N113").Value = -5957.875

$ws.Range("H126").Value = 3616.6667
$ws.Range("J126").Value = 3850
$ws.Range("L126").Value = 11550
$ws.Range("N126").Value = -16490

$ws.Range("H132").Value = 282834.94
$ws.Range("I132").Value = 410681.75
$ws.Range("J132").Value = 1571.9333
$ws.Range("K132").Value = 1232045.25
$ws.Range("L132").Value = 4715.7999
$ws.Range("M132").Value = -1229515.25
$ws.Range("N132").Value = -9775.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4297.0454
$ws.Range("I64").Value = 1933.7142
$ws.Range("J64").Value = 5399.933
$ws.Range("K64").Value = 5801.142599999999
$ws.Range("L64").Value = 16199.799
$ws.Range("M64").Value = -5531.142599999999
$ws.Range("N64").Value = -16739.799

$ws.Range("H67").Value = 4297.0454
$ws.Range("I67").Value = 1933.7142
$ws.Range("J67").Value = 5399.933
$ws.Range("K67").Value = 5801.142599999999
$ws.Range("L67").Value = 16199.799
$ws.Range("M67").Value = -4865.142599999999
$ws.Range("N67").Value = -18071.799

$ws.Range("H70").Value = 19266.666
$ws.Range("I70").Value = 19266.666
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 57799.99800000001
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = -57484.99800000001
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 19266.666
$ws.Range("I73").Value = 19266.666
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 57799.99800000001
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = -56707.99800000001
$ws.Range("M73").ClearContents()

$ws.Range("H94").Value = 4078
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4078
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = 12234
$ws.Range("N94").Value = -13586
$ws.Range("L94").ClearContents()

$ws.Range("H134").Value = 4806.92
$ws.Range("I134").Value = 2911.111
$ws.Range("J134").Value = 5873.3125
$ws.Range("K134").Value = 8733.332999999999
$ws.Range("L134").Value = 17619.9375
$ws.Range("M134").Value = -3663.332999999999
$ws.Range("N134").Value = -27759.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1692.25
$ws.Range("I113").Value = 1228.5714
$ws.Range("J113").Value = 2052.889
$ws.Range("K113").Value = 1228.5714
$ws.Range("L113").Value = 2052.889
$ws.Range("M113").Value = 941.4286
$ws.Range("N113").Value = -6392.889

$ws.Range("H126").Value = 2401.7144
$ws.Range("I126").Value = 1362.4
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 4087.2
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -1617.2
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2872.2
$ws.Range("I7").Value = 2684.9412
$ws.Range("J7").Value = 3933.3333
$ws.Range("K7").Value = 2684.9412
$ws.Range("L7").Value = 3933.3333
$ws.Range("M7").Value = -2572.9412
$ws.Range("N7").Value = -4157.3333

$ws.Range("H68").Value = 3082.2727
$ws.Range("I68").Value = 1780.4
$ws.Range("J68").Value = 4167.1665
$ws.Range("K68").Value = 1780.4
$ws.Range("L68").Value = 4167.1665
$ws.Range("M68").Value = -1031.4
$ws.Range("N68").Value = -5665.1665

$ws.Range("H71").Value = 3082.2727
$ws.Range("I71").Value = 1780.4
$ws.Range("J71").Value = 4167.1665
$ws.Range("K71").Value = 8902
$ws.Range("L71").Value = 20835.8325
$ws.Range("M71").Value = -5158
$ws.Range("N71").Value = -28323.8325

$ws.Range("H82").Value = 2024.1
$ws.Range("I82").Value = 1520.3334
$ws.Range("J82").Value = 2436.2727
$ws.Range("K82").Value = 1520.3334
$ws.Range("L82").Value = 2436.2727
$ws.Range("M82").Value = -1159.3334
$ws.Range("N82").Value = -3158.2727

$ws.Range("H85").Value = 2024.1
$ws.Range("I85").Value = 1520.3334
$ws.Range("J85").Value = 2436.2727
$ws.Range("K85").Value = 1520.3334
$ws.Range("L85").Value = 2436.2727
$ws.Range("M85").Value = -272.3334
$ws.Range("N85").Value = -4932.2727

$ws.Range("H126").Value = 2872.2
$ws.Range("I126").Value = 2684.9412
$ws.Range("J126").Value = 3933.3333
$ws.Range("K126").Value = 8054.823600000001
$ws.Range("L126").Value = 11799.9999
$ws.Range("M126").Value = -5584.823600000001
$ws.Range("N126").Value = -16739.9999

$ws.Range("H136").Value = 1497.3954
$ws.Range("I136").Value = 1228.6836
$ws.Range("J136").Value = 4530
$ws.Range("K136").Value = 3686.0508
$ws.Range("L136").Value = 13590
$ws.Range("M136").Value = -1136.0508
$ws.Range("N136").Value = -18690

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4520.6
$ws.Range("I62").Value = 3900
$ws.Range("J62").Value = 4934.3335
$ws.Range("K62").Value = 3900
$ws.Range("L62").Value = 4934.3335
$ws.Range("M62").Value = -3276
$ws.Range("N62").Value = -6182.3335

$ws.Range("H65").Value = 4520.6
$ws.Range("I65").Value = 3900
$ws.Range("J65").Value = 4934.3335
$ws.Range("K65").Value = 19500
$ws.Range("L65").Value = 24671.6675
$ws.Range("M65").Value = -16380
$ws.Range("N65").Value = -30911.6675

$ws.Range("H113").Value = 707.6667
$ws.Range("I113").Value = 569.2
$ws.Range("K113").Value = 1707.6
$ws.Range("M113").Value = 462.3999999999999
